$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (scheduled GitHub Actions update).
# D-column cells that look like plain numbers need an explicit Text
# number format first, otherwise Excel auto-coerces the assigned string
# into a numeric value instead of keeping it as text (matching the
# original inline-string cell type used for every Price cell).

$ws.Range("D2").Value = "57.577.66"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "2.446.62"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "491.61"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.59"
$ws.Range("E6").Value = "  +1.86%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").Value = "  +19.30%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.995"
$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("D9").Value = "2.476.97"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.26"
$ws.Range("E10").Value = "  +9.74%  "

$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("D14").Value = "2.877.30"
$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("D15").Value = "57.538.53"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.78"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("E17").Value = "  -2.43%  "

$ws.Range("D18").Value = "2.468.89"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.67"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "326.40"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").Value = "  +1.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.40"
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -3.32%  "

$ws.Range("E27").Value = "  -1.65%  "

$ws.Range("D28").Value = "2.563.73"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  -2.69%  "

$ws.Range("D30").Value = "0.0₃0805"
$ws.Range("E30").Value = "  -1.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.83"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.17"
$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.35"
$ws.Range("E35").Value = "  +1.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.80"
$ws.Range("E36").Value = "  +1.00%  "

$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.833"
$ws.Range("E38").Value = "  -6.15%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.32"
$ws.Range("E39").Value = "  +0.34%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.39"
$ws.Range("E40").Value = "  -1.12%  "

$ws.Range("E41").Value = "  +1.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "284.24"
$ws.Range("E42").Value = "  +7.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.101"
$ws.Range("E43").Value = "  +4.33%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.993"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.606"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0537"
$ws.Range("E46").Value = "  -3.99%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.25"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.66"
$ws.Range("E49").Value = "  -3.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.97"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").Value = "1.898.45"
$ws.Range("E51").Value = "  +3.58%  "
